$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data cell
$ws.Range("B2").Value = 100

# Add new data row
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 200
$ws.Range("C3").Value = "test"
$ws.Range("D3").Value = 31411
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("E3").Value = "Perso"
$ws.Range("F3").Value = 2000

# Convert range to an Excel Table (ListObject)
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:F3"), [Type]::Missing, 1)
$lo.Name = "Table1"
$lo.TableStyle = "TableStyleMedium9"

# Row/column sizing to mirror the final layout
$ws.Rows("1:1").AutoFit()
$ws.Columns(1).ColumnWidth = 12.666666666666666
$ws.Columns(2).ColumnWidth = 10.666666666666666
$ws.Columns(5).ColumnWidth = 13.666666666666666
$ws.Columns(6).ColumnWidth = 18.333333333333336

$ws.Range("C13").Select()
